$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("engineering")
$ws.Range("A1").Value = "Hello"
